$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '64.977.86'
$ws.Range('E2').Value = '  +3.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.383.97'
$ws.Range('E3').Value = '  +3.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.92'
$ws.Range('E5').Value = '  +3.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.64'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').Value = '  +2.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.375.59'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('E10').Value = '  +12.10%  '
$ws.Range('E11').Value = '  +4.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.02'
$ws.Range('E12').Value = '  +3.59%  '
$ws.Range('E13').Value = '  +6.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.10'
$ws.Range('E14').Value = '  +4.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.919.49'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.25'
$ws.Range('E16').Value = '  +2.92%  '
$ws.Range('E17').Value = '  +3.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.382.80'
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '64.916.88'
$ws.Range('E19').Value = '  +3.05%  '
$ws.Range('E20').Value = '  +2.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.993'
$ws.Range('E21').Value = '  +3.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '471.72'
$ws.Range('E22').Value = '  +15.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.97'
$ws.Range('E23').Value = '  +14.26%  '
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.71'
$ws.Range('E25').Value = '  +5.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.48'
$ws.Range('E26').Value = '  +0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('E27').Value = '  +7.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.77'
$ws.Range('E28').Value = '  +3.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.76'
$ws.Range('E29').Value = '  +3.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.64'
$ws.Range('E30').Value = '  +6.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.69'
$ws.Range('E31').Value = '  +7.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.49'
$ws.Range('E32').Value = '  +2.67%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '570.77'
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '61.45'
$ws.Range('E34').Value = '  +7.05%  '
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  +8.82%  '
$ws.Range('E38').Value = '  -3.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.57'
$ws.Range('E39').Value = '  +2.76%  '
$ws.Range('E40').Value = '  +2.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.369'
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.086.68'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.998'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('E44').Value = '  +4.53%  '
$ws.Range('E45').Value = '  +5.10%  '
$ws.Range('E46').Value = '  +6.34%  '
$ws.Range('E47').Value = '  +3.34%  '
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.39'
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('E51').Value = '  +4.80%  '
